$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.025.31"
$ws.Range("E2").Value = "  +5.22%  "
$ws.Range("D3").Value = "2.256.32"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "301.84"
$ws.Range("E5").Value = "  +3.23%  "
$ws.Range("D6").Value = "92.70"
$ws.Range("E6").Value = "  +5.50%  "
$ws.Range("D7").Value = "0.534"
$ws.Range("E7").Value = "  +3.91%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +2.97%  "
$ws.Range("D10").Value = "32.69"
$ws.Range("E10").Value = "  +6.76%  "
$ws.Range("D11").Value = "54.60"
$ws.Range("E11").Value = "  +9.24%  "
$ws.Range("E12").Value = "  +2.42%  "
$ws.Range("E13").Value = "  +3.46%  "
$ws.Range("D14").Value = "6.67"
$ws.Range("E14").Value = "  +3.33%  "
$ws.Range("D15").Value = "2.607.24"
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").Value = "14.14"
$ws.Range("E16").Value = "  +2.53%  "
$ws.Range("D17").Value = "2.260.67"
$ws.Range("E17").Value = "  +2.24%  "
$ws.Range("D18").Value = "0.757"
$ws.Range("E18").Value = "  +3.37%  "
$ws.Range("D19").Value = "41.913.49"
$ws.Range("E19").Value = "  +5.08%  "
$ws.Range("D20").Value = "12.16"
$ws.Range("E20").Value = "  +9.18%  "
$ws.Range("D21").Value = "0.0₃0906"
$ws.Range("E21").Value = "  +2.15%  "
$ws.Range("D22").Value = "5.95"
$ws.Range("E22").Value = "  +3.67%  "
$ws.Range("D23").Value = "67.26"
$ws.Range("E23").Value = "  +2.52%  "
$ws.Range("D24").Value = "242.17"
$ws.Range("E24").Value = "  +2.13%  "
$ws.Range("D25").Value = "2.58"
$ws.Range("E25").Value = "  +4.92%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +4.27%  "
$ws.Range("D28").Value = "23.95"
$ws.Range("E28").Value = "  +3.39%  "
$ws.Range("E29").Value = "  +4.57%  "
$ws.Range("E30").Value = "  +2.12%  "
$ws.Range("D31").Value = "34.13"
$ws.Range("E31").Value = "  +6.96%  "
$ws.Range("D32").Value = "158.62"
$ws.Range("E32").Value = "  +1.01%  "
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").Value = "5.15"
$ws.Range("E34").Value = "  +3.71%  "
$ws.Range("D35").Value = "0.0744"
$ws.Range("E35").Value = "  +4.41%  "
$ws.Range("D36").Value = "3.05"
$ws.Range("E36").Value = "  +3.36%  "
$ws.Range("E37").Value = "  +2.55%  "
$ws.Range("E38").Value = "  +5.73%  "
$ws.Range("D39").Value = "16.63"
$ws.Range("E39").Value = "  +8.48%  "
$ws.Range("E40").Value = "  +4.02%  "
$ws.Range("D41").Value = "1.80"
$ws.Range("E41").Value = "  +4.36%  "
$ws.Range("E42").Value = "  +5.38%  "
$ws.Range("D43").Value = "2.050.59"
$ws.Range("E43").Value = "  -2.78%  "
$ws.Range("D44").Value = "19.87"
$ws.Range("E44").Value = "  +11.26%  "
$ws.Range("D45").Value = "0.0279"
$ws.Range("E45").Value = "  +3.56%  "
$ws.Range("D46").Value = "10.12"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("E47").Value = "  +8.14%  "
$ws.Range("D48").Value = "2.02"
$ws.Range("E48").Value = "  -2.72%  "
$ws.Range("D49").Value = "2.478.17"
$ws.Range("E49").Value = "  +1.88%  "
$ws.Range("E50").Value = "  +2.86%  "
$ws.Range("E51").Value = "  +4.44%  "
